$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs.Item(2)

$metaXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Chilli Pop, a colorful Mexican cuisine-inspired online slot game. Play for free and enjoy free spins, bonus features, and multipliers.</w:t></w:r></w:p>
"@
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------
# 2) At the bottom of the document: remove the duplicated bold title
#    paragraph ("Play Chilli Pop for Free...") and replace the text of
#    the remaining italic paragraph with the new image-prompt text,
#    keeping its italic formatting.
# ---------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldTitlePara = $d.Paragraphs.Item($count - 1)
$boldTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastRange = $lastPara.Range
$lastRange.MoveEnd(1, -1) | Out-Null
$lastRange.Text = "Create a vibrant feature image for Chilli Pop that showcases the game's Mexican theme and fun characters. The image should be in a cartoon style and feature a happy Maya warrior with glasses as the main focus. Surround the warrior with symbols from the game, such as tomatoes, garlic, peppers, and onions that have been transformed into wacky characters. Use bright colors to make the image pop and ensure that it captures the spirit of the game's exciting gameplay and cluster-based winning combinations. Add the Chilli Pop logo to the center of the image to tie it all together and make it clear which game it represents."

Write-Host "Edit applied."
